$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill B1:Q1 with the same header text as B1 ("Tasa con 30kV")
$headerTasa = $ws.Range("B1").Text
$ws.Range("B1:Q1").Value = $headerTasa

# Fill B2:Q2 with the same header text as B2 ("R(30kV)/Imp/s")
$headerR = $ws.Range("B2").Text
$ws.Range("B2:Q2").Value = $headerR

# Update row 4 values
$ws.Range("A4").Value = 30
$ws.Range("B4").Value = 13183
$ws.Range("C4:Q4").Value = 0

# Update row 5 values
$ws.Range("A5").Value = 40
$ws.Range("B5").Value = 12352
$ws.Range("C5:Q5").Value = 0

# Remove old rows 6-9 (data that is no longer present)
$ws.Range("A6:Q9").Delete()

# Update the selected range to match the new data extent
$ws.Range("A1:Q5").Select()
